$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.139.40'
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').Value = '2.927.02'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.89'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.21'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.46%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  +0.46%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.99'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +4.64%  '
$ws.Range('E10').Value = '  -0.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.441'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000225'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.126'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').Value = '3.414.10'
$ws.Range('E15').Value = '  +0.92%  '
$ws.Range('D16').Value = '61.116.35'
$ws.Range('E16').Value = '  +0.70%  '
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('D18').Value = '2.931.98'
$ws.Range('E18').Value = '  +1.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '435.94'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.79%  '
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('E21').Value = '  -0.97%  '
$ws.Range('E22').Value = '  +0.57%  '
$ws.Range('E23').Value = '  +0.55%  '
$ws.Range('E24').Value = '  +2.93%  '
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.88'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.99'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.69%  '
$ws.Range('E31').Value = '  +4.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.71'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.02%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').Value = '0.0₃0866'
$ws.Range('E34').Value = '  +1.88%  '
$ws.Range('E35').Value = '  -0.28%  '
$ws.Range('E36').Value = '  +0.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.99'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.31%  '
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('E39').Value = '  +0.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.60'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '41.97'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +3.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.287'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '376.73'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.09%  '
$ws.Range('E44').Value = '  -0.60%  '
$ws.Range('D45').Value = '2.688.48'
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '133.63'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.99'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.25%  '
$ws.Range('E49').Value = '  -0.42%  '
$ws.Range('E50').Value = '  -1.66%  '
$ws.Range('E51').Value = '  -0.07%  '
